$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 33

$ws.Cells.Item($newRow, 1).Value = 10002
$ws.Cells.Item($newRow, 2).Value = 110032
$ws.Cells.Item($newRow, 3).Value = 10032
$ws.Cells.Item($newRow, 4).Value = "eng"
$ws.Cells.Item($newRow, 5).Value = $true
$ws.Cells.Item($newRow, 6).Value = "superadmin"
$ws.Cells.Item($newRow, 7).Value = "now()"

$ws.Range("C30").Select()
